$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Alternativa [ID Inválido] (passo 5)" / "Regressa a 3" cells and
# replace with the revised "Exceção" wording, clearing the old "Regressa a 3" cell.
$ws.Range("B15").Value = "Exceção 1 [ID Inválido] (passo 5)"
$ws.Range("D16").ClearContents()

# Update the selection / view state to match the saved workbook
$ws.Range("E16").Select()
$excel.ActiveWindow.ScrollRow = 2
